# Regenerate s_val data to filter save games: update computed stat values
# for rows 2-4 (columns B-G) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 21.98653043760045)
    3 = @(0.01253208636536152, 0.3048912486333797, 3.223369029078222, 2797.565817734744, 2801.106610098821)
    4 = @(0.1169995834814548, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.67637130870356)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
